{"js": "// Remove the \"Full reproducibility\" heading paragraph and the paragraph\n// that follows it (\"All the code chunks together can be found in this\n// gist.\") \u2014 these two paragraphs are deleted entirely from the document.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst toDelete = [];\nfor (const para of paragraphs.items) {\n  const t = (para.text || \"\").trim();\n  if (t === \"Full reproducibility\" || t.indexOf(\"All the code chunks together can be found in\") !== -1) {\n    toDelete.push(para);\n  }\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Full reproducibility\" heading paragraph and the paragraph\n# that follows it (\"All the code chunks together can be found in this\n# gist.\") \u2014 these two paragraphs are deleted entirely from the document.\n\n$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $para = $d.Paragraphs.Item($i)\n    $text = $para.Range.Text.Trim()\n    if ($text -eq \"Full reproducibility\" -or $text.StartsWith(\"All the code chunks together can be found in\")) {\n        $para.Range.Delete()\n    }\n}\n"}
